# Search Listing page: add one more test step row to the "Test Steps" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")

# Duplicate the last existing row (row 9) into the new row 10 so that the
# new row inherits the same cell style (fill/border) as the rest of the
# table, then overwrite its contents with the new step's data.
$ws.Range("A9:E9").Copy($ws.Range("A10:E10"))

$ws.Range("A10").Value = "searchListingPageTestCases"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "addProduct"
$ws.Range("D10").Value = "search_list_projectNames|search_list_addBtn"
$ws.Range("E10").Value = "Apple - Royal Gala"

# Match the saved selection/view state recorded for the sheet after the edit.
$ws.Range("C13").Select()
